$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.737.62"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.605.90"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.40%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "213.28"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  +0.47%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "28.20"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +5.55%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.254"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("E10").Value = "  +0.57%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0911"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "1.834.99"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "1.611.34"
$ws.Range("E13").Value = "  +0.33%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.550"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +4.02%  "
$ws.Range("D15").Value = "29.706.94"
$ws.Range("E15").Value = "  +0.09%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.77"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.17%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "64.15"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.69%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "242.07"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("E19").Value = "  +3.41%  "
$ws.Range("D20").Value = "0.0₃0699"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("E22").Value = "  -0.60%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "9.42"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("E24").Value = "  -0.62%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "155.11"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "15.48"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.72%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.109"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("E31").Value = "  +0.92%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.25"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").Value = "1.426.57"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("E35").Value = "  +3.31%  "
$ws.Range("E36").Value = "  +2.43%  "
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  +1.53%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.547"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +2.10%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "56.60"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.49%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.0493"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +5.25%  "
$ws.Range("E43").Value = "  +2.01%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.95"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.985"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +17.69%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "66.32"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("D49").Value = "1.743.37"
$ws.Range("E49").Value = "  -0.08%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "86.62"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").Value = "0.0₆0104"
$ws.Range("E51").Value = "  +1.31%  "
